# 5.2.1.1a.xlsx — add a new "2023" data column (column Q) to the table,
# mirroring the layout/styling already used for the 2022 column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (the "Число женщин, обратившихся по факту..." header row) grows a
# touch taller to fit the extra column of wrapped text.
$ws.Rows.Item(5).RowHeight = 27

# New column Q values for 2023, one per data row.
$ws.Cells.Item(4, 17).Value = 2023   # header year
$ws.Cells.Item(6, 17).Value = 1209
$ws.Cells.Item(7, 17).Value = "-"
$ws.Cells.Item(8, 17).Value = 373
$ws.Cells.Item(9, 17).Value = 115
$ws.Cells.Item(10, 17).Value = 781

# Column Q (rows 3-10) should carry the same cell formatting as column P,
# the previous last column in the table (borders, number formats, etc.).
$ws.Range("P3:P10").Copy()
$ws.Range("Q3:Q10").PasteSpecial(-4122) # xlPasteFormats

# Clear the lingering "P7" selection left over from editing, settling back
# on the sheet's natural top-left cell.
$ws.Range("A1").Select()
